$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data updates
$ws.Range("A2").Value = 968754
$ws.Range("B2").Value = "TestUser112"

# F2: becomes a text email address styled/linked like F3 and F4
$ws.Range("F3").Copy($ws.Range("F2"))
$wb.Worksheets.Item(1).Hyperlinks.Add($ws.Range("F2"), "mailto:test@d")
$ws.Range("F2").Style = $ws.Range("F3").Style

# New column width for column B
$ws.Columns("B:B").ColumnWidth = 13.5

# Update the active selection shown in the sheet view
$ws.Range("H2").Select() | Out-Null
